# Applies the "Remove Fn.boxes, allow only Fn.seqnum and Fn.box ... rename Fn.box -> Fn.string"
# edit to the "Protocol" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# New row 8: add the "seqnum, nextSeqNum" message row (previously blank), with Y/Y/Y/Y
# across the C2S HTTP / C2S Socket-like / S2C HTTP / S2C Socket-like columns.
$ws.Range("A8").Value = "seqnum, nextSeqNum"
$ws.Range("B8").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("D8").Value = "Y"
$ws.Range("E8").Value = "Y"

# Row 9: rename Fn.boxes -> Fn.box (renders as "box, string"), and update the
# comment to reflect that a single string (a previous "box") is sent, rather
# than "boxes".
$ws.Range("A9").Value = "box, string"
$ws.Range("F9").Value = "Both parties send strings (previous boxes) - this is the whole point of Minerva."

# Update the active selection to reflect where the author's cursor ended up
# after adding the new row.
$ws.Activate() | Out-Null
$ws.Range("A8").Select() | Out-Null
